$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 135
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = 34
$ws.Range("N4").Value = -528
$ws.Range("H6").Value = 797.75
$ws.Range("I6").Value = 797.75
$ws.Range("K6").Value = 2393.25
$ws.Range("M6").Value = -2281.25
$ws.Range("H8").Value = 252
$ws.Range("I8").Value = 252
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 756
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -617
$ws.Range("N8").ClearContents()
$ws.Range("H40").Value = 1975
$ws.Range("I40").Value = 1950
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1950
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1775
$ws.Range("N40").Value = -2350
$ws.Range("H51").Value = 2900.875
$ws.Range("J51").Value = 2921.2
$ws.Range("L51").Value = 2921.2
$ws.Range("N51").Value = -3889.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 352.66666
$ws.Range("I5").Value = 466.23077
$ws.Range("K5").Value = 466.23077
$ws.Range("M5").Value = -354.23077
$ws.Range("H102").Value = 1298.8182
$ws.Range("I102").Value = 898.5714
$ws.Range("J102").Value = 1999.25
$ws.Range("K102").Value = 898.5714
$ws.Range("L102").Value = 1999.25
$ws.Range("M102").Value = 723.4286
$ws.Range("N102").Value = -5243.25
$ws.Range("H132").Value = 9519
$ws.Range("I132").Value = 7779.8335
$ws.Range("J132").Value = 12364.909
$ws.Range("K132").Value = 23339.5005
$ws.Range("L132").Value = 37094.727
$ws.Range("M132").Value = -20809.5005
$ws.Range("N132").Value = -42154.727

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 352.66666
$ws.Range("I4").Value = 466.23077
$ws.Range("K4").Value = 466.23077
$ws.Range("M4").Value = -351.23077
$ws.Range("H105").Value = 3370
$ws.Range("I105").Value = 3555
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3555
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1808
$ws.Range("N105").Value = -6494
$ws.Range("H107").Value = 1258.2084
$ws.Range("I107").Value = 892.61536
$ws.Range("J107").Value = 1690.2727
$ws.Range("K107").Value = 892.61536
$ws.Range("L107").Value = 1690.2727
$ws.Range("M107").Value = 1027.38464
$ws.Range("N107").Value = -5530.2727

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.48649
$ws.Range("J7").Value = 76.84999999999999
$ws.Range("L7").Value = 76.84999999999999
$ws.Range("N7").Value = -302.85
$ws.Range("H105").Value = 1485.8
$ws.Range("I105").Value = 1443
$ws.Range("J105").Value = 1550
$ws.Range("K105").Value = 1443
$ws.Range("L105").Value = 1550
$ws.Range("M105").Value = 304
$ws.Range("N105").Value = -5044
$ws.Range("H107").Value = 749.6
$ws.Range("I107").Value = 749.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 749.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1170.4
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 1189.2858
$ws.Range("I122").Value = 770.3333
$ws.Range("J122").Value = 1503.5
$ws.Range("K122").Value = 2310.9999
$ws.Range("L122").Value = 4510.5
$ws.Range("M122").Value = 139.0001000000002
$ws.Range("N122").Value = -9410.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 675619.25
$ws.Range("I113").Value = 492.5
$ws.Range("J113").Value = 875656.75
$ws.Range("K113").Value = 1477.5
$ws.Range("L113").Value = 2626970.25
$ws.Range("M113").Value = 692.5
$ws.Range("N113").Value = -2631310.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 801042.4
$ws.Range("I14").Value = 801042.4
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 801042.4
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -800874.4
$ws.Range("N14").ClearContents()
$ws.Range("H132").Value = 7716.8696
$ws.Range("I132").Value = 10223.571
$ws.Range("J132").Value = 3817.5557
$ws.Range("K132").Value = 30670.713
$ws.Range("L132").Value = 11452.6671
$ws.Range("M132").Value = -28140.713
$ws.Range("N132").Value = -16512.6671

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 800
$ws.Range("J18").Value = 800
$ws.Range("L18").Value = 800
$ws.Range("N18").Value = -1144
$ws.Range("H46").Value = 1702.2142
$ws.Range("I46").Value = 1892.3334
$ws.Range("J46").Value = 1360
$ws.Range("K46").Value = 1892.3334
$ws.Range("L46").Value = 1360
$ws.Range("M46").Value = -1704.3334
$ws.Range("N46").Value = -1736
$ws.Range("H55").Value = 222.35294
$ws.Range("I55").Value = 205
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 205
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -32
$ws.Range("N55").Value = -846
$ws.Range("H82").Value = 1300.8
$ws.Range("I82").Value = 1251
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 1251
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -890
$ws.Range("N82").Value = -2222
$ws.Range("H85").Value = 1300.8
$ws.Range("I85").Value = 1251
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 1251
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -3
$ws.Range("N85").Value = -3996
$ws.Range("H132").Value = 7323.207
$ws.Range("I132").Value = 10087.4
$ws.Range("J132").Value = 4361.5713
$ws.Range("K132").Value = 30262.2
$ws.Range("L132").Value = 13084.7139
$ws.Range("M132").Value = -27732.2
$ws.Range("N132").Value = -18144.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 25300.5
$ws.Range("J103").Value = 25300.5
$ws.Range("L103").Value = 25300.5
$ws.Range("N103").Value = -27644.5
$ws.Range("H107").Value = 420
$ws.Range("I107").Value = 375
$ws.Range("J107").Value = 492
$ws.Range("K107").Value = 1125
$ws.Range("L107").Value = 1476
$ws.Range("M107").Value = 795
$ws.Range("N107").Value = -5316
